$wb = $excel.ActiveWorkbook

# ===== Sheet: 展览 =====
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("22:23").Insert()
$ws1.Range("A22").Style = $ws1.Range("A21").Style
$ws1.Range("A22").Value = 21
$ws1.Range("B22").NumberFormat = "@"
$ws1.Range("B22").Value = "2024-04-20"
$ws1.Range("C22").Value = "广州·SISP动漫游戏嘉年华之地下城探险（免费活动）"
$ws1.Range("D22").Value = "西湾路150号 悦汇城"
$ws1.Range("E22").Value = "2024.04.20 13:00-04.21 19:00"
$ws1.Range("F22").Value = 4
$ws1.Range("G22").Value = 48
$ws1.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=83210"
$ws1.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202403/TZO1ioLk1711079685372.jpeg"

$ws1.Range("A23").Style = $ws1.Range("A22").Style
$ws1.Range("A23").Value = 22
$ws1.Range("B23").NumberFormat = "@"
$ws1.Range("B23").Value = "2024-04-20"
$ws1.Range("C23").Value = "广州·次元方舟动漫游戏嘉年华"
$ws1.Range("D23").Value = "东沙大道16号 广州国际医药港"
$ws1.Range("E23").Value = "2024.04.20 09:00-04.21 18:00"
$ws1.Range("F23").Value = 19
$ws1.Range("G23").Value = 65
$ws1.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=83217"
$ws1.Range("I23").Value = "//i0.hdslb.com/bfs/openplatform/202403/wGYGykoT1710753036838.jpeg"

$ws1.Range("29:29").Insert()
$ws1.Range("A29").Style = $ws1.Range("A28").Style
$ws1.Range("A29").Value = 28
$ws1.Range("B29").NumberFormat = "@"
$ws1.Range("B29").Value = "2024-05-05"
$ws1.Range("C29").Value = "广州·AI动漫展4.0"
$ws1.Range("D29").Value = "奥体南路12号 优托邦(奥体旗舰店)"
$ws1.Range("E29").Value = "2024.05.05 10:00-05.05 17:00"
$ws1.Range("F29").Value = 1
$ws1.Range("G29").Value = 55
$ws1.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=83221"
$ws1.Range("I29").Value = "//i1.hdslb.com/bfs/openplatform/202403/G3zFqk0D1710904038569.jpeg"

$ws1.Range("35:36").Insert()
$ws1.Range("A35").Style = $ws1.Range("A34").Style
$ws1.Range("A35").Value = 34
$ws1.Range("B35").NumberFormat = "@"
$ws1.Range("B35").Value = "2024-05-19"
$ws1.Range("C35").Value = "广州·蓝锁only3.0"
$ws1.Range("D35").Value = "大石街石北工业大道644号 巨大创意产业园"
$ws1.Range("E35").Value = "2024.05.19 10:00-05.19 17:00"
$ws1.Range("F35").Value = 0
$ws1.Range("G35").Value = 65
$ws1.Range("H35").Value = "https://show.bilibili.com/platform/detail.html?id=83226"
$ws1.Range("I35").Value = "//i2.hdslb.com/bfs/openplatform/202403/sdu8lzow1711033421343.jpeg"

$ws1.Range("A36").Style = $ws1.Range("A35").Style
$ws1.Range("A36").Value = 35
$ws1.Range("B36").NumberFormat = "@"
$ws1.Range("B36").Value = "2024-05-25"
$ws1.Range("C36").Value = "广州·EVAonly海边集市同人展"
$ws1.Range("D36").Value = "黄边三横路一街1号 设计殿堂"
$ws1.Range("E36").Value = "2024.05.25 10:00-05.26 17:00"
$ws1.Range("F36").Value = 0
$ws1.Range("G36").Value = 78
$ws1.Range("H36").Value = "https://show.bilibili.com/platform/detail.html?id=82918"
$ws1.Range("I36").Value = "//i1.hdslb.com/bfs/openplatform/202403/VtRyhhHg1710487512346.png"

# ---- F column value updates for 展览 ----
$ws1.Range("F2").Value = 622
$ws1.Range("F3").Value = 287
$ws1.Range("F5").Value = 766
$ws1.Range("F6").Value = 431
$ws1.Range("F8").Value = 195
$ws1.Range("F10").Value = 267
$ws1.Range("F11").Value = 7070
$ws1.Range("F14").Value = 95
$ws1.Range("F16").Value = 564
$ws1.Range("F17").Value = 386
$ws1.Range("F20").Value = 20
$ws1.Range("F24").Value = 81
$ws1.Range("F25").Value = 3
$ws1.Range("F26").Value = 194
$ws1.Range("F28").Value = 342
$ws1.Range("F30").Value = 1050
$ws1.Range("F32").Value = 59
$ws1.Range("F33").Value = 2026
$ws1.Range("F34").Value = 555
$ws1.Range("F38").Value = 544
$ws1.Range("F39").Value = 7

# ===== Sheet: 演出 =====
$ws2 = $wb.Worksheets.Item("演出")

# ---- F column value updates for 演出 ----
$ws2.Range("F4").Value = 57
$ws2.Range("F10").Value = 137

# ===== Sheet: 本地生活 =====
$ws3 = $wb.Worksheets.Item("本地生活")

# ---- F column value updates for 本地生活 ----
$ws3.Range("F2").Value = 339

# ===== Sheet: 全部类型 =====
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("29:30").Insert()
$ws4.Range("A29").Style = $ws4.Range("A28").Style
$ws4.Range("A29").Value = 28
$ws4.Range("B29").NumberFormat = "@"
$ws4.Range("B29").Value = "2024-04-20"
$ws4.Range("C29").Value = "广州·SISP动漫游戏嘉年华之地下城探险（免费活动）"
$ws4.Range("D29").Value = "西湾路150号 悦汇城"
$ws4.Range("E29").Value = "2024.04.20 13:00-04.21 19:00"
$ws4.Range("F29").Value = 4
$ws4.Range("G29").Value = 48
$ws4.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=83210"
$ws4.Range("I29").Value = "//i1.hdslb.com/bfs/openplatform/202403/TZO1ioLk1711079685372.jpeg"

$ws4.Range("A30").Style = $ws4.Range("A29").Style
$ws4.Range("A30").Value = 29
$ws4.Range("B30").NumberFormat = "@"
$ws4.Range("B30").Value = "2024-04-20"
$ws4.Range("C30").Value = "广州·次元方舟动漫游戏嘉年华"
$ws4.Range("D30").Value = "东沙大道16号 广州国际医药港"
$ws4.Range("E30").Value = "2024.04.20 09:00-04.21 18:00"
$ws4.Range("F30").Value = 19
$ws4.Range("G30").Value = 65
$ws4.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=83217"
$ws4.Range("I30").Value = "//i0.hdslb.com/bfs/openplatform/202403/wGYGykoT1710753036838.jpeg"

$ws4.Range("39:39").Insert()
$ws4.Range("A39").Style = $ws4.Range("A38").Style
$ws4.Range("A39").Value = 38
$ws4.Range("B39").NumberFormat = "@"
$ws4.Range("B39").Value = "2024-05-05"
$ws4.Range("C39").Value = "广州·AI动漫展4.0"
$ws4.Range("D39").Value = "奥体南路12号 优托邦(奥体旗舰店)"
$ws4.Range("E39").Value = "2024.05.05 10:00-05.05 17:00"
$ws4.Range("F39").Value = 1
$ws4.Range("G39").Value = 55
$ws4.Range("H39").Value = "https://show.bilibili.com/platform/detail.html?id=83221"
$ws4.Range("I39").Value = "//i1.hdslb.com/bfs/openplatform/202403/G3zFqk0D1710904038569.jpeg"

$ws4.Range("45:46").Insert()
$ws4.Range("A45").Style = $ws4.Range("A44").Style
$ws4.Range("A45").Value = 44
$ws4.Range("B45").NumberFormat = "@"
$ws4.Range("B45").Value = "2024-05-19"
$ws4.Range("C45").Value = "广州·蓝锁only3.0"
$ws4.Range("D45").Value = "大石街石北工业大道644号 巨大创意产业园"
$ws4.Range("E45").Value = "2024.05.19 10:00-05.19 17:00"
$ws4.Range("F45").Value = 0
$ws4.Range("G45").Value = 65
$ws4.Range("H45").Value = "https://show.bilibili.com/platform/detail.html?id=83226"
$ws4.Range("I45").Value = "//i2.hdslb.com/bfs/openplatform/202403/sdu8lzow1711033421343.jpeg"

$ws4.Range("A46").Style = $ws4.Range("A45").Style
$ws4.Range("A46").Value = 45
$ws4.Range("B46").NumberFormat = "@"
$ws4.Range("B46").Value = "2024-05-25"
$ws4.Range("C46").Value = "广州·EVAonly海边集市同人展"
$ws4.Range("D46").Value = "黄边三横路一街1号 设计殿堂"
$ws4.Range("E46").Value = "2024.05.25 10:00-05.26 17:00"
$ws4.Range("F46").Value = 0
$ws4.Range("G46").Value = 78
$ws4.Range("H46").Value = "https://show.bilibili.com/platform/detail.html?id=82918"
$ws4.Range("I46").Value = "//i1.hdslb.com/bfs/openplatform/202403/VtRyhhHg1710487512346.png"

# ---- F column value updates for 全部类型 ----
$ws4.Range("F2").Value = 339
$ws4.Range("F3").Value = 622
$ws4.Range("F4").Value = 287
$ws4.Range("F6").Value = 766
$ws4.Range("F8").Value = 431
$ws4.Range("F10").Value = 195
$ws4.Range("F12").Value = 267
$ws4.Range("F13").Value = 7070
$ws4.Range("F17").Value = 95
$ws4.Range("F19").Value = 564
$ws4.Range("F20").Value = 386
$ws4.Range("F22").Value = 57
$ws4.Range("F24").Value = 20
$ws4.Range("F31").Value = 81
$ws4.Range("F33").Value = 3
$ws4.Range("F36").Value = 194
$ws4.Range("F38").Value = 342
$ws4.Range("F40").Value = 1050
$ws4.Range("F42").Value = 59
$ws4.Range("F43").Value = 2026
$ws4.Range("F44").Value = 555
$ws4.Range("F48").Value = 544
$ws4.Range("F49").Value = 7
